# "Import notes when we have them" — add a "Notes" header column (H1)
# to the import template, matching the yellow-fill header style used by
# the other non-bold header cells (A1, D1, E1, F1, G1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell with its text.
$ws.Range("H1").Value = "Notes"

# Match formatting of the existing plain header cells (yellow fill,
# regular weight font) rather than the bold Origin/Destination style.
$ws.Range("H1").Interior.Color = $ws.Range("A1").Interior.Color
$ws.Range("H1").Font.Bold = $ws.Range("A1").Font.Bold

# Leave the selection on the newly added header cell.
[void]$ws.Range("H1").Select()
